# Updated cryptos list on Sat Aug 12 05:39:03 UTC 2023 with GitHub Actions
# Refreshes Price (D) and Volume(1h) (E) columns on the crypto sheet; rows 39/40 (MXToken/VeChain) swapped order.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.389.19"
$ws.Range("E2").Value = "  +0.05%  "
$ws.Range("D3").Value = "1.846.32"
$ws.Range("E3").Value = "  +0.00%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9994"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.73"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.18%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6296"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.94%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07588"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.44%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2929"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.25%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.47"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.71%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07733"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.03%  "
$ws.Range("D12").Value = "1.842.28"
$ws.Range("E12").Value = "  +0.01%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.00001091"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +9.72%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.000"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.18%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6782"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.73%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "83.64"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.67%  "
$ws.Range("D17").Value = "2.093.08"
$ws.Range("E17").Value = "  -7.56%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.152"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.22%  "
$ws.Range("D19").Value = "29.408.43"
$ws.Range("E19").Value = "  +0.04%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "228.56"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.78%  "
$ws.Range("E21").Value = "  -0.08%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.425"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.73%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.001"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "157.39"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.29%  "
$ws.Range("E26").Value = "  -0.50%  "
$ws.Range("E27").Value = "  -0.14%  "
$ws.Range("E28").Value = "  -0.25%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.462"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.300"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.24%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.05606"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.73%  "
$ws.Range("E32").Value = "  -0.56%  "
$ws.Range("E33").Value = "  +0.07%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.843"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.14%  "
$ws.Range("E35").Value = "  -0.02%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7085"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.10%  "
$ws.Range("E37").Value = "  -0.35%  "
$ws.Range("D38").Value = "1.229.57"
$ws.Range("E38").Value = "  -2.05%  "
$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.768"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.79%  "
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01796"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.62%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.447"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.00%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9059"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.23%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9998"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.07%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.88"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.12%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "66.02"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.66%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000122"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.25%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.174"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.57%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4017"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.16%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.991"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.67%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.678"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.52%  "
$ws.Range("E51").Value = "  -0.65%  "
